$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 229421
$ws.Range("C3").Value = 4471
$ws.Range("C4").Value = 308
$ws.Range("C5").Value = 7238
$ws.Range("C6").Value = 166282
$ws.Range("C7").Value = 303824
$ws.Range("C8").Value = 162309
$ws.Range("C9").Value = 203684
$ws.Range("C10").Value = 185942
$ws.Range("C11").Value = 47
